$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns hold plain text values (e.g. "27.971.13",
# "  -0.88%  ") rather than numbers, so force Text formatting before assigning,
# then restore the original (default) style so no formatting changes leak in.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.971.13"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "1.869.10"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "312.64"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "0.5035"
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("D8").Value = "0.3824"
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("D9").Value = "0.08952"
$ws.Range("E9").Value = "  -7.13%  "
$ws.Range("D10").Value = "1.116"
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").Value = "41.57"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").Value = "6.361"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("D14").Value = "1.870.89"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "7.232"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "1.0000"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "91.06"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").Value = "0.06646"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "18.19"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "6.118"
$ws.Range("E22").Value = "  -1.90%  "
$ws.Range("D23").Value = "28.002.45"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").Value = "11.48"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").Value = "2.262"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("D26").Value = "2.084.26"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").Value = "2.493"
$ws.Range("E27").Value = "  -6.43%  "
$ws.Range("D28").Value = "157.52"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").Value = "20.71"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").Value = "126.23"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").Value = "0.1065"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  -3.60%  "
$ws.Range("D33").Value = "5.598"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").Value = "3.604"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("D35").Value = "9.425"
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("D36").Value = "0.06588"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").Value = "0.02395"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").Value = "0.2184"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").Value = "1.204"
$ws.Range("E40").Value = "  -3.09%  "
$ws.Range("D41").Value = "0.6374"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "11.46"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "4.893"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "0.6009"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "13.13"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("D47").Value = "1.279"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "3.658"
$ws.Range("E48").Value = "  -2.45%  "
$ws.Range("D49").Value = "1.233"
$ws.Range("E49").Value = "  +3.78%  "
$ws.Range("D50").Value = "1.994"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").Value = "120.72"
$ws.Range("E51").Value = "  -2.14%  "

$dataRange.Style = "Normal"
